$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Update the two Avg_Time_ms values that changed (D2 and D3)
$ws.Range("D2").Value = 0.6102165500000001
$ws.Range("D3").Value = 1.2739587

# Make sure dependent objects (e.g. the chart cache) are refreshed
$excel.CalculateFullRebuild()
$wb.RefreshAll()
